# Embedding layer - ablation study
# Adds two new rows of ablation results ("without position embeddings" and
# "without token type embeddings") to the Abstractive block (rows 7-8) and
# to the Extractive block (rows 19-20) of the "Sheet2" worksheet (3rd tab,
# internally named "Sheet2" even though the 2nd tab is confusingly named
# "Results").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows into the Abstractive results block -----------
# (pushes the old rows 7-16 down to 9-18)
$ws.Rows("7:8").Insert(-4121)   # xlShiftDown

# --- Insert the two new rows into the Extractive results block ------------
# (these row numbers already account for the shift above; pushes the old
# rows 17-28 down to 21-32)
$ws.Rows("19:20").Insert(-4121) # xlShiftDown

# --- Populate the Extractive-block new rows first, so the new shared ------
# --- strings are allocated in "position" -> "token type" order, matching --
# --- the target workbook's shared string table.                       -----
$ws.Range("F19").Value = "without position embeddings"
$ws.Range("J19").Value = 25.1
$ws.Range("K19").Value = 7.54
$ws.Range("L19").Value = 19

$ws.Range("F20").Value = "without token type embeddings"
$ws.Range("J20").Value = 31.9
$ws.Range("K20").Value = 12.64
$ws.Range("L20").Value = 23.9

# --- Populate the Abstractive-block new rows -------------------------------
$ws.Range("F7").Value = "without token type embeddings"
$ws.Range("J7").Value = 47.2
$ws.Range("K7").Value = 24.77
$ws.Range("L7").Value = 37.5

$ws.Range("F8").Value = "without position embeddings"
$ws.Range("J8").Value = 36
$ws.Range("K8").Value = 12.85
$ws.Range("L8").Value = 27.51

# --- Restore the active cell the author ended on (J36) ---------------------
[void]$ws.Activate()
[void]$ws.Range("J36").Select()
